$d = $word.ActiveDocument

# First paragraph of the document ("**ID__AFFARS_MP5305_3__ID** ")
$para = $d.Paragraphs.Item(1)

# Remove the trailing space run, leaving only the **ID__...__ID** run.
$pEnd = $para.Range.End
$trailingSpace = $d.Range($pEnd - 2, $pEnd - 1)
$trailingSpace.Delete()

# Update the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$para.Format.LeftIndent = 11.25

# Add a paragraph border (box) with 5 twips of spacing on every edge.
$para.Format.Borders.DistanceFromTop = 5
$para.Format.Borders.DistanceFromLeft = 5
$para.Format.Borders.DistanceFromBottom = 5
$para.Format.Borders.DistanceFromRight = 5
